# Update RMSE table header units and refresh the active cell selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Unit update: RMSE output is in units of model output, not squared units.
$ws.Range("B1").Value = "DO RMSE (mg/L)"
$ws.Range("C1").Value = "DOC RMSE (mg/L)"

# Update the selected/active cell on the sheet (was C9, now C2).
$ws.Activate()
$ws.Range("C2").Select()
